$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FilesTab query text in cell B4: remove the "File Type" and
# "Breed" columns from the RETURN clause (corrected ICDC Breed script).
$newFilesTabQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed  IN ['Saint Bernard']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n         coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n          coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Cells.Item(4, 2).Value = $newFilesTabQuery

# Row 4 height shrinks because the text now has fewer lines.
$ws.Rows.Item(4).RowHeight = 217.5

# Update the view: select B4 and scroll so row 4 is at the top.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
